$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 103
$ws.Range("I11").Value = 103
$ws.Range("K11").Value = 103
$ws.Range("M11").Value = 37

# Row 40
$ws.Range("H40").Value = 4986.533
$ws.Range("I40").Value = 3219.9
$ws.Range("J40").Value = 8519.799999999999
$ws.Range("K40").Value = 3219.9
$ws.Range("L40").Value = 8519.799999999999
$ws.Range("M40").Value = -3044.9
$ws.Range("N40").Value = -8869.799999999999

# Row 43
$ws.Range("H43").Value = 1500
$ws.Range("J43").Value = 1500
$ws.Range("L43").Value = 1500
$ws.Range("N43").Value = -1638

# Row 87
$ws.Range("H87").Value = 66434.664
$ws.Range("J87").Value = 66434.664
$ws.Range("L87").Value = 66434.664
$ws.Range("N87").Value = -68930.664

# Row 90
$ws.Range("H90").Value = 66434.664
$ws.Range("J90").Value = 66434.664
$ws.Range("L90").Value = 199303.992
$ws.Range("N90").Value = -211783.992

# Row 138
$ws.Range("H138").Value = 5589.189
$ws.Range("J138").Value = 5736.364
$ws.Range("L138").Value = 17209.092
$ws.Range("N138").Value = -27489.092


$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1972.7142
$ws.Range("I2").Value = 1972.7142
$ws.Range("K2").Value = 1972.7142
$ws.Range("M2").Value = -1859.7142

# Row 32
$ws.Range("H32").Value = 15032.412
$ws.Range("I32").Value = 12862.5625
$ws.Range("K32").Value = 12862.5625
$ws.Range("M32").Value = -12575.5625

# Row 88
$ws.Range("H88").Value = 2749.1667
$ws.Range("I88").Value = 1065
$ws.Range("J88").Value = 4433.3335
$ws.Range("K88").Value = 1065
$ws.Range("L88").Value = 4433.3335
$ws.Range("M88").Value = -659
$ws.Range("N88").Value = -5245.3335

# Row 91
$ws.Range("H91").Value = 2749.1667
$ws.Range("I91").Value = 1065
$ws.Range("J91").Value = 4433.3335
$ws.Range("K91").Value = 1065
$ws.Range("L91").Value = 4433.3335
$ws.Range("M91").Value = 339
$ws.Range("N91").Value = -7241.3335

# Row 116
$ws.Range("H116").Value = 1972.7142
$ws.Range("I116").Value = 1972.7142
$ws.Range("K116").Value = 1972.7142
$ws.Range("M116").Value = 321.2858000000001


$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1972.7142
$ws.Range("I3").Value = 1972.7142
$ws.Range("K3").Value = 1972.7142
$ws.Range("M3").Value = -1858.7142

# Row 9
$ws.Range("H9").Value = 58999.8
$ws.Range("J9").Value = 58999.8
$ws.Range("L9").Value = 58999.8
$ws.Range("N9").Value = -59335.8

# Row 130
$ws.Range("H130").Value = 98466.164
$ws.Range("J130").Value = 98466.164
$ws.Range("L130").Value = 98466.164
$ws.Range("N130").Value = -108506.164

# Row 131
$ws.Range("H131").Value = 1000000
$ws.Range("J131").Value = 1000000
$ws.Range("L131").Value = 1000000
$ws.Range("N131").Value = -1010080


$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 210.5
$ws.Range("I7").Value = 159.25
$ws.Range("J7").Value = 261.75
$ws.Range("K7").Value = 159.25
$ws.Range("L7").Value = 261.75
$ws.Range("M7").Value = -46.25
$ws.Range("N7").Value = -487.75

# Row 99
$ws.Range("H99").Value = 3567.8
$ws.Range("I99").Value = 3634.3333
$ws.Range("K99").Value = 3634.3333
$ws.Range("M99").Value = -2136.3333

# Row 126
$ws.Range("H126").Value = 3567.8
$ws.Range("I126").Value = 3634.3333
$ws.Range("K126").Value = 10902.9999
$ws.Range("M126").Value = -8432.999899999999

# Row 132
$ws.Range("H132").Value = 3487
$ws.Range("I132").Value = 2938.4
$ws.Range("K132").Value = 8815.200000000001
$ws.Range("M132").Value = -6285.200000000001

# Row 134
$ws.Range("H134").Value = 2704.75
$ws.Range("I134").Value = 2000
$ws.Range("J134").Value = 2805.4285
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 8416.2855
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -13486.2855


$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 110
$ws.Range("J7").Value = 99.5
$ws.Range("L7").Value = 298.5
$ws.Range("N7").Value = -522.5

# Row 60
$ws.Range("H60").Value = 924.25
$ws.Range("I60").Value = 899
$ws.Range("K60").Value = 2697
$ws.Range("M60").Value = -2446

# Row 68
$ws.Range("H68").Value = 1999.8
$ws.Range("J68").Value = 1999.8
$ws.Range("L68").Value = 5999.4
$ws.Range("N68").Value = -7621.4

# Row 71
$ws.Range("H71").Value = 1999.8
$ws.Range("J71").Value = 1999.8
$ws.Range("L71").Value = 17998.2
$ws.Range("N71").Value = -26110.2

# Row 92
$ws.Range("H92").Value = 380
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

# Row 98
$ws.Range("H98").Value = 2561.5
$ws.Range("I98").Value = 2875.8
$ws.Range("J98").Value = 990
$ws.Range("K98").Value = 8627.400000000001
$ws.Range("L98").Value = 2970
$ws.Range("M98").Value = -7129.400000000001
$ws.Range("N98").Value = -5966

# Row 109
$ws.Range("H109").Value = 4273.88
$ws.Range("J109").Value = 5000
$ws.Range("L109").Value = 15000
$ws.Range("N109").Value = -17080

# Row 113
$ws.Range("H113").Value = 1248.5
$ws.Range("J113").Value = 1166.6666
$ws.Range("L113").Value = 3499.9998
$ws.Range("N113").Value = -7839.9998

# Row 132
$ws.Range("H132").Value = 10040
$ws.Range("I132").Value = 1300
$ws.Range("K132").Value = 11700
$ws.Range("M132").Value = -9170

# Row 134
$ws.Range("H134").Value = 2131.75
$ws.Range("I134").Value = 2131.75
$ws.Range("K134").Value = 6395.25
$ws.Range("M134").Value = -1325.25


$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2452.7
$ws.Range("I68").Value = 2744.4119
$ws.Range("J68").Value = 799.6667
$ws.Range("K68").Value = 2744.4119
$ws.Range("L68").Value = 799.6667
$ws.Range("M68").Value = -1995.4119
$ws.Range("N68").Value = -2297.6667

# Row 71
$ws.Range("H71").Value = 2452.7
$ws.Range("I71").Value = 2744.4119
$ws.Range("J71").Value = 799.6667
$ws.Range("K71").Value = 13722.0595
$ws.Range("L71").Value = 3998.3335
$ws.Range("M71").Value = -9978.059499999999
$ws.Range("N71").Value = -11486.3335


$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2470.182
$ws.Range("I81").Value = 2470.182
$ws.Range("K81").Value = 4940.364
$ws.Range("M81").Value = -3879.364

# Row 84
$ws.Range("H84").Value = 2470.182
$ws.Range("I84").Value = 2470.182
$ws.Range("K84").Value = 24701.82
$ws.Range("M84").Value = -19397.82

# Row 113
$ws.Range("H113").Value = 768.75
$ws.Range("I113").Value = 807.1429000000001
$ws.Range("K113").Value = 2421.4287
$ws.Range("M113").Value = -251.4287000000004

# Row 126
$ws.Range("H126").Value = 1614.5
$ws.Range("I126").Value = 1487.25
$ws.Range("J126").Value = 1639.95
$ws.Range("K126").Value = 4461.75
$ws.Range("L126").Value = 4919.85
$ws.Range("M126").Value = -1991.75
$ws.Range("N126").Value = -9859.85

# Row 132
$ws.Range("H132").Value = 4290
$ws.Range("I132").Value = 1682.7142
$ws.Range("K132").Value = 5048.142599999999
$ws.Range("M132").Value = -2518.142599999999

